$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header labels "Std" / "Relative std" (bold, matching D6/E6, D9/E9, D12/E12 style)
$ws.Range("D15").Value = "Std"
$ws.Range("E15").Value = "Relative std"
$ws.Range("D15:E15").Font.Bold = $true

# Add new formulas for standard deviation and relative standard deviation
$ws.Range("D16").Formula = "=STDEV(B2:B31)"
$ws.Range("E16").Formula = "= (D16 / E4) * 100"

# Update the active selection to E16
$ws.Range("E16").Select()
